$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Each D-column value here is a dotted/numeric-looking piece of text
# (e.g. "29.181.68", "1.000"). Excel auto-detects plain numeric literals
# and would silently convert them to numbers (dropping text like trailing
# zeros or multi-dot grouping), so each cell is temporarily forced to Text
# format, written, then restored to the workbook default ("Normal") style
# so no stray formatting is left behind.

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '29.181.68'
$cell.Style = "Normal"

$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.826.21'
$cell.Style = "Normal"

$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '0.9993'
$cell.Style = "Normal"

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '236.63'
$cell.Style = "Normal"

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.6123'
$cell.Style = "Normal"

$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"

$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.07103'
$cell.Style = "Normal"

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.2816'
$cell.Style = "Normal"

$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '23.54'
$cell.Style = "Normal"

$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.07665'
$cell.Style = "Normal"

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '1.824.80'
$cell.Style = "Normal"

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.817'
$cell.Style = "Normal"

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '0.00001007'
$cell.Style = "Normal"

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.6322'
$cell.Style = "Normal"

$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '2.067.10'
$cell.Style = "Normal"

$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '78.95'
$cell.Style = "Normal"

$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '5.870'
$cell.Style = "Normal"

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '29.174.43'
$cell.Style = "Normal"

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '227.60'
$cell.Style = "Normal"

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '11.79'
$cell.Style = "Normal"

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.Style = "Normal"

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '7.005'
$cell.Style = "Normal"

$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '154.99'
$cell.Style = "Normal"

$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '0.1317'
$cell.Style = "Normal"

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '8.038'
$cell.Style = "Normal"

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '16.61'
$cell.Style = "Normal"

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.489'
$cell.Style = "Normal"

$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '0.06345'
$cell.Style = "Normal"

$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '1.453'
$cell.Style = "Normal"

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.821'
$cell.Style = "Normal"

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.795'
$cell.Style = "Normal"

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.748'
$cell.Style = "Normal"

$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.6495'
$cell.Style = "Normal"

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '2.544'
$cell.Style = "Normal"

$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '2.748'
$cell.Style = "Normal"

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.216.39'
$cell.Style = "Normal"

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '6.586'
$cell.Style = "Normal"

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.01742'
$cell.Style = "Normal"

$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.9249'
$cell.Style = "Normal"

$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '0.9993'
$cell.Style = "Normal"

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '101.23'
$cell.Style = "Normal"

$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '1.973.92'
$cell.Style = "Normal"

$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '63.00'
$cell.Style = "Normal"

$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.624'
$cell.Style = "Normal"

$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '8.630'
$cell.Style = "Normal"

$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.4562'
$cell.Style = "Normal"

$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.05518'
$cell.Style = "Normal"

# --- Volume(1h) (column E) updates ---
# These are already padded percentage strings ("  -0.55%  "), which Excel
# keeps as text because of the surrounding whitespace, so a plain Value
# assignment is sufficient.
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('E6').Value = '  -2.58%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -4.48%  '
$ws.Range('E9').Value = '  -2.59%  '
$ws.Range('E10').Value = '  -5.57%  '
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('E15').Value = '  -6.39%  '
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('E18').Value = '  -5.84%  '
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  -5.26%  '
$ws.Range('E28').Value = '  -4.53%  '
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('E30').Value = '  -10.39%  '
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('E32').Value = '  -5.66%  '
$ws.Range('E33').Value = '  -5.88%  '
$ws.Range('E34').Value = '  -1.01%  '
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('E36').Value = '  -6.74%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('E40').Value = '  -3.22%  '
$ws.Range('E41').Value = '  -5.25%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('E46').Value = '  -3.44%  '
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('E48').Value = '  -4.79%  '
$ws.Range('E49').Value = '  -3.32%  '
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('E51').Value = '  -2.70%  '
